$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.159.90'
$ws.Range('E2').Value = '  +1.94%  '
$ws.Range('D3').Value = '2.379.07'
$ws.Range('E3').Value = '  +0.71%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.695'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +7.29%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '243.56'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +4.43%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '77.54'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +7.82%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.609'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +27.74%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.105'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +6.94%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '57.84'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +1.83%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '32.62'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +20.26%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '7.55'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +20.38%  '
$ws.Range('E14').Value = '  +2.47%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '17.29'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +8.35%  '
$ws.Range('B16').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C16').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D16').Value = '2.730.46'
$ws.Range('E16').Value = '  +0.32%  '
$ws.Range('E17').Value = '  +8.40%  '
$ws.Range('D18').Value = '2.374.72'
$ws.Range('E18').Value = '  +0.70%  '
$ws.Range('D19').Value = '44.443.88'
$ws.Range('E19').Value = '  +2.79%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0000105'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +2.15%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.74'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +6.40%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '78.92'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +6.26%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '258.38'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +3.61%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.57'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +4.96%  '
$ws.Range('E26').Value = '  +0.80%  '
$ws.Range('B27').Value = 'ImmutableX'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.81'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +20.01%  '
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.98'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +10.05%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '23.16'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +3.34%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.31'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +2.20%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '175.73'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.17%  '
$ws.Range('E32').Value = '  +0.97%  '
$ws.Range('E33').Value = '  +7.33%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.40'
$ws.Range('D34').Style = "Normal"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0764'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +10.78%  '
$ws.Range('E36').Value = '  +6.58%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.92'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +6.74%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.51'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +2.52%  '
$ws.Range('E39').Value = '  +0.72%  '
$ws.Range('E40').Value = '  +9.87%  '
$ws.Range('E41').Value = '  +3.24%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '19.20'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +4.02%  '
$ws.Range('E43').Value = '  -0.03%  '
$ws.Range('E44').Value = '  +18.53%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.23'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +4.36%  '
$ws.Range('E46').Value = '  +6.20%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.57'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +15.12%  '
$ws.Range('E48').Value = '  +6.37%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '103.21'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +3.33%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '4.50'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.35%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '54.85'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +9.11%  '
